$p = $ppt.ActivePresentation

# Remove the "Sample Heading With Picture" sample slide (slide 10 of 11)
# leaving "THANK YOU, THE END" as the new final slide.
$p.Slides.Item(10).Delete()
